$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.019935090384098
$ws.Cells.Item(2, 4).Value = 1.025772014519371
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.031010098666089
$ws.Cells.Item(2, 9).Value = 1.029236789521426
$ws.Cells.Item(2, 10).Value = 1.025135044189006
$ws.Cells.Item(2, 11).Value = 1.028596608461961
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.033819454449791
$ws.Cells.Item(2, 14).Value = 1.026590853536389

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.020888860600917
$ws.Cells.Item(3, 4).Value = 1.026471517666629
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.032179677360702
$ws.Cells.Item(3, 9).Value = 1.029395159313322
$ws.Cells.Item(3, 10).Value = 1.02572576208916
$ws.Cells.Item(3, 11).Value = 1.02910377790706
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.034796538839376
$ws.Cells.Item(3, 14).Value = 1.027182410323716

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.021506131702722
$ws.Cells.Item(4, 4).Value = 1.026924044155253
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.032936790795867
$ws.Cells.Item(4, 9).Value = 1.02949625073165
$ws.Cells.Item(4, 10).Value = 1.026107534605203
$ws.Cells.Item(4, 11).Value = 1.029431179792321
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.035428526798164
$ws.Cells.Item(4, 14).Value = 1.027564725000522

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.021765660360305
$ws.Cells.Item(5, 4).Value = 1.027114261640713
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.033255157681296
$ws.Cells.Item(5, 9).Value = 1.029538417788735
$ws.Cells.Item(5, 10).Value = 1.026267920797977
$ws.Cells.Item(5, 11).Value = 1.029568634404598
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.035694154584416
$ws.Cells.Item(5, 14).Value = 1.027725338960086

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.02180923798869
$ws.Cells.Item(6, 4).Value = 1.0271461985191
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.033308617371369
$ws.Cells.Item(6, 9).Value = 1.029545478357807
$ws.Cells.Item(6, 10).Value = 1.026294843823499
$ws.Cells.Item(6, 11).Value = 1.029591702763389
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.035738751146864
$ws.Cells.Item(6, 14).Value = 1.027752300219393

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.02150959942822
$ws.Cells.Item(7, 4).Value = 1.026926585948715
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.032941044529189
$ws.Cells.Item(7, 9).Value = 1.029496815474006
$ws.Cells.Item(7, 10).Value = 1.026109678130034
$ws.Cells.Item(7, 11).Value = 1.029433017195864
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.035432076365049
$ws.Cells.Item(7, 14).Value = 1.027566871569404

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.020257396845266
$ws.Cells.Item(8, 4).Value = 1.026008434749557
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.031405297644351
$ws.Cells.Item(8, 9).Value = 1.029290597784223
$ws.Cells.Item(8, 10).Value = 1.025334775320225
$ws.Cells.Item(8, 11).Value = 1.0287681679348
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.034149716740197
$ws.Cells.Item(8, 14).Value = 1.026790868308722

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.018051775185107
$ws.Cells.Item(9, 4).Value = 1.024389821639464
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.028701535045164
$ws.Cells.Item(9, 9).Value = 1.028916627718677
$ws.Cells.Item(9, 10).Value = 1.023965784925101
$ws.Cells.Item(9, 11).Value = 1.027590744355421
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.031888106625069
$ws.Cells.Item(9, 14).Value = 1.025419933790212

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016581997332842
$ws.Cells.Item(10, 4).Value = 1.023310323809645
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.026900630164901
$ws.Cells.Item(10, 9).Value = 1.028660214869633
$ws.Cells.Item(10, 10).Value = 1.023050790540383
$ws.Cells.Item(10, 11).Value = 1.026801882328875
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.030379057885528
$ws.Cells.Item(10, 14).Value = 1.024503640008516

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015945721378117
$ws.Cells.Item(11, 4).Value = 1.022842800760564
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.026121192633271
$ws.Cells.Item(11, 9).Value = 1.028547505223868
$ws.Cells.Item(11, 10).Value = 1.022654039542685
$ws.Cells.Item(11, 11).Value = 1.026459375519791
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.029725309115907
$ws.Cells.Item(11, 14).Value = 1.024106325578894

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.015709402477307
$ws.Cells.Item(12, 4).Value = 1.022669128740722
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.025831729195445
$ws.Cells.Item(12, 9).Value = 1.028505387529081
$ws.Cells.Item(12, 10).Value = 1.022506585925843
$ws.Cells.Item(12, 11).Value = 1.026332014893296
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.029482429107272
$ws.Cells.Item(12, 14).Value = 1.023958662561002

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.015760092671327
$ws.Cells.Item(13, 4).Value = 1.022706382559117
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.025893817635658
$ws.Cells.Item(13, 9).Value = 1.028514433331056
$ws.Cells.Item(13, 10).Value = 1.022538218960665
$ws.Cells.Item(13, 11).Value = 1.026359340432124
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.029534529899998
$ws.Cells.Item(13, 14).Value = 1.023990340518362

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.01592618672788
$ws.Cells.Item(14, 4).Value = 1.02282844524024
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.026097264389885
$ws.Cells.Item(14, 9).Value = 1.028544028907668
$ws.Cells.Item(14, 10).Value = 1.02264185266429
$ws.Cells.Item(14, 11).Value = 1.026448850669892
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.029705233572655
$ws.Cells.Item(14, 14).Value = 1.024094121393733

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.016028525719925
$ws.Cells.Item(15, 4).Value = 1.022903650371465
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.026222621834461
$ws.Cells.Item(15, 9).Value = 1.028562230293102
$ws.Cells.Item(15, 10).Value = 1.022705693870767
$ws.Cells.Item(15, 11).Value = 1.026503982564877
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.029810403285672
$ws.Cells.Item(15, 14).Value = 1.024158053262046

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016624227994062
$ws.Cells.Item(16, 4).Value = 1.023341349883126
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.026952366526338
$ws.Cells.Item(16, 9).Value = 1.028667659639573
$ws.Cells.Item(16, 10).Value = 1.023077109989907
$ws.Cells.Item(16, 11).Value = 1.026824593947734
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.030422438245078
$ws.Cells.Item(16, 14).Value = 1.024529996834677

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.01699793590305
$ws.Cells.Item(17, 4).Value = 1.023615882851956
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.027410213668432
$ws.Cells.Item(17, 9).Value = 1.02873334276864
$ws.Cells.Item(17, 10).Value = 1.023309941872359
$ws.Cells.Item(17, 11).Value = 1.027025457895555
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.030806265442993
$ws.Cells.Item(17, 14).Value = 1.024763159365109

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.017215927630457
$ws.Cells.Item(18, 4).Value = 1.02377600422051
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.02767730353321
$ws.Cells.Item(18, 9).Value = 1.028771492398271
$ws.Cells.Item(18, 10).Value = 1.023445695448437
$ws.Cells.Item(18, 11).Value = 1.027142529216912
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.031030114438535
$ws.Cells.Item(18, 14).Value = 1.024899105726836

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.017290259584618
$ws.Cells.Item(19, 4).Value = 1.023830599904323
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.02776838029443
$ws.Cells.Item(19, 9).Value = 1.028784472911279
$ws.Cells.Item(19, 10).Value = 1.023491974854436
$ws.Cells.Item(19, 11).Value = 1.027182432362655
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.03110643593043
$ws.Cells.Item(19, 14).Value = 1.0249454508549

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.016957839108426
$ws.Cells.Item(20, 4).Value = 1.023586428993997
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.027361087348813
$ws.Cells.Item(20, 9).Value = 1.028726312368594
$ws.Cells.Item(20, 10).Value = 1.023284966727526
$ws.Cells.Item(20, 11).Value = 1.027003916319613
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.030765087617781
$ws.Cells.Item(20, 14).Value = 1.024738148752705

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.015877275550683
$ws.Cells.Item(21, 4).Value = 1.022792501178711
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.02603735289596
$ws.Cells.Item(21, 9).Value = 1.028535320710971
$ws.Cells.Item(21, 10).Value = 1.022611337399712
$ws.Cells.Item(21, 11).Value = 1.026422495950715
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.029654966940943
$ws.Cells.Item(21, 14).Value = 1.02406356279398

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.015198011560537
$ws.Cells.Item(22, 4).Value = 1.022293251643993
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.025205382732105
$ws.Cells.Item(22, 9).Value = 1.028413776996263
$ws.Cells.Item(22, 10).Value = 1.022187322004886
$ws.Cells.Item(22, 11).Value = 1.026056133710784
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.028956708573876
$ws.Cells.Item(22, 14).Value = 1.023638945248659

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.015558089956044
$ws.Cells.Item(23, 4).Value = 1.022557920140958
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.025646396287909
$ws.Cells.Item(23, 9).Value = 1.028478347891084
$ws.Cells.Item(23, 10).Value = 1.022412145768047
$ws.Cells.Item(23, 11).Value = 1.026250424938832
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.029326895364343
$ws.Cells.Item(23, 14).Value = 1.02386408828735

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016975957079909
$ws.Cells.Item(24, 4).Value = 1.023599737952642
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.027383285308451
$ws.Cells.Item(24, 9).Value = 1.028729489605007
$ws.Cells.Item(24, 10).Value = 1.023296252085153
$ws.Cells.Item(24, 11).Value = 1.027013650306209
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.030783694199805
$ws.Cells.Item(24, 14).Value = 1.024749450136835

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.018621870406218
$ws.Cells.Item(25, 4).Value = 1.024808350554223
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.029400238150842
$ws.Cells.Item(25, 9).Value = 1.029014560694303
$ws.Cells.Item(25, 10).Value = 1.024320115173833
$ws.Cells.Item(25, 11).Value = 1.027895828926742
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.032473016835119
$ws.Cells.Item(25, 14).Value = 1.025774767228539
